$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header fields ---
$ws.Range("B2").Value = "Richard Dobson"
$ws.Range("F2").Value = "Week 8"

# --- Activity rows (4-8): Week 8, Mon-Fri 16-20 Sep 2019, 9am-1pm, alternating Group/Individual ---
$activity = "Work on csv_output function "

$ws.Range("A4").Value = $activity
$ws.Range("C4").Value = "G"
$ws.Range("D4").Value = 43724
$ws.Range("E4").Value = 0.375
$ws.Range("F4").Value = 0.54166666666666663
$ws.Range("G4").Value = 4

$ws.Range("A5").Value = $activity
$ws.Range("C5").Value = "I"
$ws.Range("D5").Value = 43725
$ws.Range("E5").Value = 0.375
$ws.Range("F5").Value = 0.54166666666666663
$ws.Range("H5").Value = 4

$ws.Range("A6").Value = $activity
$ws.Range("C6").Value = "G"
$ws.Range("D6").Value = 43726
$ws.Range("E6").Value = 0.375
$ws.Range("F6").Value = 0.54166666666666663
$ws.Range("G6").Value = 4

$ws.Range("A7").Value = $activity
$ws.Range("C7").Value = "I"
$ws.Range("D7").Value = 43727
$ws.Range("E7").Value = 0.375
$ws.Range("F7").Value = 0.54166666666666663
$ws.Range("H7").Value = 4

$ws.Range("A8").Value = $activity
$ws.Range("C8").Value = "G"
$ws.Range("D8").Value = 43728
$ws.Range("E8").Value = 0.375
$ws.Range("F8").Value = 0.54166666666666663
$ws.Range("G8").Value = 4

# --- Visually separate each day's entry with a double top border on rows 5-8 (row 4 already had it) ---
$ws.Range("A5:D5").Borders.Item(8).LineStyle = -4119
$ws.Range("A6:D6").Borders.Item(8).LineStyle = -4119
$ws.Range("A7:D7").Borders.Item(8).LineStyle = -4119
$ws.Range("A8:D8").Borders.Item(8).LineStyle = -4119

# --- Column widths to fit the new data ---
$ws.Range("B1").ColumnWidth = 19
$ws.Range("D1").ColumnWidth = 13

# --- Selection matches the author's last saved position ---
$ws.Range("M14:N14").Select()
